$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 new rows before row 198, shifting existing rows 198-268 down to 201-271.
$xlShiftDown = [Microsoft.Office.Interop.Excel.XlInsertShiftDirection]::xlShiftDown
$ws.Range("A198:T200").Insert($xlShiftDown)

# Fill the 3 newly inserted rows with the new data block (date 2021-11-10 / serial 44510).
$newRows = @(
    @{ Row = 198; Quality = "Especial"; Vol = 400; PMin = 12500; PMax = 13000; PProm = 12750; PKg = 1821 },
    @{ Row = 199; Quality = "Primera";  Vol = 360; PMin = 10500; PMax = 11000; PProm = 10750; PKg = 1536 },
    @{ Row = 200; Quality = "Segunda";  Vol = 300; PMin = 8500;  PMax = 9000;  PProm = 8750;  PKg = 1250 }
)

foreach ($r in $newRows) {
    $row = $r.Row
    $ws.Cells.Item($row, 1).Value  = 2
    $ws.Cells.Item($row, 2).Value  = "Comercializadora del Agro de Limarí"
    $ws.Cells.Item($row, 3).Value  = "Coquimbo"
    $ws.Cells.Item($row, 4).Value  = 44510
    $ws.Cells.Item($row, 5).Value  = 4
    $ws.Cells.Item($row, 6).Value  = "Fruta"
    $ws.Cells.Item($row, 7).Value  = 100101
    $ws.Cells.Item($row, 8).Value  = "Berries"
    $ws.Cells.Item($row, 9).Value  = 100112025
    $ws.Cells.Item($row, 10).Value = "Frutilla"
    $ws.Cells.Item($row, 11).Value = "Sin especificar"
    $ws.Cells.Item($row, 12).Value = $r.Quality
    $ws.Cells.Item($row, 13).Value = $r.Vol
    $ws.Cells.Item($row, 14).Value = $r.PMin
    $ws.Cells.Item($row, 15).Value = $r.PMax
    $ws.Cells.Item($row, 16).Value = $r.PProm
    $ws.Cells.Item($row, 17).Value = "$/bandeja 7 kilos"
    $ws.Cells.Item($row, 18).Value = "Provincia de Melipilla"
    $ws.Cells.Item($row, 19).Value = $r.PKg
    $ws.Cells.Item($row, 20).Value = 7
}

Write-Host ("Final dims: {0}" -f $ws.UsedRange.Address())
